# Adds a new date column (column 126 = "DV") to the attendance tracker,
# the next tracked date after the previous last column (DU =
# 2026-02-20, serial 46073), dated 2026-02-23 (serial 46076), and fills
# in each player's attendance mark for that date, mirroring the pattern
# already present in column DU.
#
# Note: the header cell (row 1, no formulas depend on it) gets its
# style copied straight from DU1 via Range.Copy so it reuses the
# existing date-number-format style instead of minting a new one.
# The data rows (2-32) instead get their style set by nudging the
# (already-applied) center alignment, which Excel folds back onto the
# existing shared style too -- Copy() on those rows was found to leave
# the row's COUNTA/COUNTIF formulas stuck at their old cached values,
# so plain .Value assignments are used there to keep auto-recalc
# working normally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$duCol = 125  # column DU (previous last date column)
$dvCol = 126  # column DV (new date column)

# Header date cell DV1: copy DU1 (style + value), then overwrite with
# the new date serial number (2026-02-23).
$ws.Cells.Item(1, $duCol).Copy($ws.Cells.Item(1, $dvCol))
$ws.Cells.Item(1, $dvCol).Value = 46076

# Attendance letter codes for each player row in the new DV column.
# Rows 12, 21 and 23 are players whose tracked range already ends
# before column DU, so they get no new cell at all. Rows 16, 17 and 25
# get an empty (but present) cell, like their DU neighbour.
$values = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "B"
    7  = "P"
    8  = "B"
    9  = "P"
    10 = "P"
    11 = "P"
    13 = "B"
    14 = "P"
    15 = "P"
    16 = ""
    17 = ""
    18 = "P"
    19 = "P"
    20 = "P"
    22 = "P"
    24 = "P"
    25 = ""
    26 = "P"
    27 = "P"
    28 = "A"
    29 = "P"
    30 = "P"
    31 = "P"
    32 = "P"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Cells.Item($row, $dvCol)
    $cell.Value = $values[$row]
    $cell.HorizontalAlignment = -4108   # xlCenter -- matches existing style "4"
}

# Update the frozen pane / active selection to reflect the new last column.
$ws.Application.ActiveWindow.ScrollColumn = $dvCol
$ws.Range("DV30").Select()
